# Calibration data clean-up: the curvature samples were recorded out of
# order; re-sort the data rows (A2:D18) in ascending order of the
# "time (s)" column (column A) while keeping each row's B/C/D curvature
# values together with its original time value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = $ws.Range("A2:D18").Value()
$rowCount = $data.GetLength(0)

# Pull each worksheet row into its own 4-element array so it can be
# reordered as a unit.
$rowList = New-Object System.Collections.ArrayList
for ($i = 1; $i -le $rowCount; $i++) {
    $r = @($data[$i, 1], $data[$i, 2], $data[$i, 3], $data[$i, 4])
    [void]$rowList.Add($r)
}

# Sort the rows by the time column (column A).
$sortedRows = $rowList | Sort-Object { $_[0] }

# Write the reordered rows back into the same range.
for ($i = 0; $i -lt $sortedRows.Count; $i++) {
    $targetRow = $i + 2
    $row = $sortedRows[$i]
    $ws.Cells.Item($targetRow, 1).Value = $row[0]
    $ws.Cells.Item($targetRow, 2).Value = $row[1]
    $ws.Cells.Item($targetRow, 3).Value = $row[2]
    $ws.Cells.Item($targetRow, 4).Value = $row[3]
}
